{"js": "// Append two new paragraphs at the end of the document body:\n//   1. An empty paragraph that only carries italic paragraph-mark formatting\n//      (mirrors the italic run that precedes it).\n//   2. A regular paragraph with the new closing text about the missing pencils.\nconst body = context.document.body;\n\nconst newText =\n  \"The next thing they know as they enter classroom again is that their \" +\n  \"pencils are gone! No one stole them but they are literally lost! \" +\n  \"Nowhere to be found. And then he remembers a thought that was well too \" +\n  \"known to all in the school life, \\u201CAakha jhimik, maal gayab!\\u201D \" +\n  \"The teacher for the next class enters now and they guffaw thinking \" +\n  \"about the poor soul who stole the pencil.\";\n\n// Use a flat-OPC OOXML payload so the inserted markup matches exactly\n// (an empty paragraph carrying only the italic paragraph mark, followed by\n// a plain paragraph holding the new sentence) without Office.js silently\n// materializing an inherited-formatting run in the empty paragraph.\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:pStyle w:val=\"NoSpacing\"/>' +\n  '<w:jc w:val=\"both\"/>' +\n  '<w:rPr>' +\n  '<w:i/>' +\n  '<w:iCs/>' +\n  '</w:rPr>' +\n  '</w:pPr>' +\n  '</w:p>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:pStyle w:val=\"NoSpacing\"/>' +\n  '<w:jc w:val=\"both\"/>' +\n  '</w:pPr>' +\n  '<w:r>' +\n  '<w:t>' + newText + '</w:t>' +\n  '</w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nconst endRange = body.getRange(\"End\");\nendRange.insertOoxml(flatOpc, \"End\");\nawait context.sync();\n", "ps1": "# Append two new paragraphs at the end of the document body:\n#   1. An empty paragraph that only carries italic paragraph-mark formatting\n#      (mirrors the italic run that precedes it).\n#   2. A regular paragraph with the new closing text about the missing pencils.\n$d = $word.ActiveDocument\n\n$openQuote = [char]0x201C\n$closeQuote = [char]0x201D\n$newText = \"The next thing they know as they enter classroom again is that their pencils are gone! No one stole them but they are literally lost! Nowhere to be found. And then he remembers a thought that was well too known to all in the school life, $($openQuote)Aakha jhimik, maal gayab!$($closeQuote) The teacher for the next class enters now and they guffaw thinking about the poor soul who stole the pencil.\"\n\n$xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"NoSpacing\"/><w:jc w:val=\"both\"/><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr></w:p><w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:pStyle w:val=\"NoSpacing\"/><w:jc w:val=\"both\"/></w:pPr><w:r><w:t>' + $newText + '</w:t></w:r></w:p>'\n\n$rng = $d.Content\n$rng.Collapse(0)\n$rng.InsertXML($xml)\n"}
